$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update title / last-updated timestamp text (row 1)
$ws.Range("A1").Value = "Datos actualizados a 4 de Abril de 2020 a las 15:20"

# Update province/city labels in column A whose shared-string slot now
# points at a different label (rows whose number stayed the same but
# whose associated text shuffled around).
$ws.Range("A13").Value = "Albacete"
$ws.Range("A14").Value = "Araba/Alava"
$ws.Range("A15").Value = "La Rioja"
$ws.Range("A17").Value = "Toledo"
$ws.Range("A18").Value = "A Coruña"
$ws.Range("A19").Value = "Malaga"
$ws.Range("A46").Value = "Cuenca"
$ws.Range("A47").Value = "Gran Canaria"
$ws.Range("A48").Value = "Huesca"

# Update the numeric data (Casos totales, Casos activos, Recuperados, Muertes)
$ws.Range("B9").Value = 3496
$ws.Range("C9").Value = 657
$ws.Range("D9").Value = 7678
$ws.Range("E9").Value = 312

$ws.Range("B13").Value = 2548
$ws.Range("C13").Value = 657
$ws.Range("D13").Value = 7678
$ws.Range("E13").Value = 194

$ws.Range("B14").Value = 2539
$ws.Range("C14").Value = 3098
$ws.Range("D14").Value = 4612
$ws.Range("E14").Value = 176

$ws.Range("B15").Value = 2405
$ws.Range("C15").Value = 843
$ws.Range("D15").Value = 1434
$ws.Range("E15").Value = 128

$ws.Range("B17").Value = 1994
$ws.Range("C17").Value = 657
$ws.Range("D17").Value = 7678
$ws.Range("E17").Value = 276

$ws.Range("B18").Value = 1969
$ws.Range("C18").Value = 333
$ws.Range("D18").Value = 1788
$ws.Range("E18").Value = 67

$ws.Range("B19").Value = 1905
$ws.Range("C19").Value = 162
$ws.Range("D19").Value = 1633
$ws.Range("E19").Value = 110

$ws.Range("B37").Value = 837
$ws.Range("C37").Value = 657
$ws.Range("D37").Value = 7678
$ws.Range("E37").Value = 118

$ws.Range("B46").Value = 449
$ws.Range("C46").Value = 657
$ws.Range("D46").Value = 7678
$ws.Range("E46").Value = 89

$ws.Range("B47").Value = 435
$ws.Range("C47").Value = 123
$ws.Range("D47").Value = 1564
$ws.Range("E47").Value = 21

$ws.Range("B48").Value = 393
$ws.Range("C48").Value = 44
$ws.Range("D48").Value = 319
$ws.Range("E48").Value = 30
